# Adding the delay wwhile executing scripts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "test run" values used by the automation (NewClient / NewProject
# test markers) and add a trailing space to the Wms project value.
$ws.Range("E2").Value = "SFGTest4091"
$ws.Range("G2").Value = "SFGTest4091"
$ws.Range("H2").Value = "Wms "

# Move the active selection back to G2 (was left on G6 from a previous run).
[void]$ws.Range("G2").Select()

# Add the dropdown data validations used while filling out the CreateJobRequest form.
$ws.Range("D2:D4").Validation.Add(3, 1, 1, '"India,Australia,Canada"')
$ws.Range("C2:C6").Validation.Add(3, 1, 1, '"IT,BPO,PST,SSS-Shared Services,SHILOH,GC-IT,DIGITAL"')
$ws.Range("F2:F5").Validation.Add(3, 1, 1, '"APL Logistics,Arbor Health,Arke"')
$ws.Range("I2:I5").Validation.Add(3, 1, 1, '".Net,Ab Initio,Admin"')
$ws.Range("J2:J4").Validation.Add(3, 1, 1, '"Associate,Architect,Associate Technical Architect"')
$ws.Range("K2:K4").Validation.Add(3, 1, 1, '"Analyst - Systems,Account Manager,Account Manager - APR"')
$ws.Range("P2:P3").Validation.Add(3, 1, 1, '"Billed,Unbilled"')
$ws.Range("R2:R3").Validation.Add(3, 1, 1, '"Yes,No"')
$ws.Range("S2:S3").Validation.Add(3, 1, 1, '"Confirmed,Pipeline"')
$ws.Range("Q2:Q3").Validation.Add(3, 1, 1, '"Addtional Billing,Replacement"')
$ws.Range("X2:X4").Validation.Add(3, 1, 1, '"Fulltime,Contract,Contract to Hire"')
$ws.Range("AI2:AI3").Validation.Add(3, 1, 1, '"Desktop,Laptop"')
